# Gantt chart.xlsx update
# Replaces the generic "Activity NN" placeholder rows with the real
# project task list (Planning / Software Design Document / Implementation /
# Controlling / Closing) together with updated plan/actual periods, renames
# the sheet title, updates the highlighted-cell selection and print scale.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data: Activity name, Plan Start, Plan Duration, Actual Start,
#     Actual Duration, Percent Complete ---------------------------------
$rows = @(
    @{ R = 5;  B = "1.0 Planning";                     C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 6;  B = "1.1 Project plan";                 C = 3;  D = 2;  E = 3;  F = 2;  G = 1 }
    @{ R = 7;  B = "1.2 Scope";                        C = 3;  D = 2;  E = 3;  F = 2;  G = 1 }
    @{ R = 8;  B = "1.3 WBS";                           C = 3;  D = 2;  E = 3;  F = 2;  G = 1 }
    @{ R = 9;  B = "1.4 Activity definition";           C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 10; B = "1.5 Scheduling and Gantt chart";    C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 11; B = "2.0 Software Design Document";      C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 12; B = "  2.1 System Vision";                C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 13; B = "2.2 Formal requirements";            C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 14; B = "2.3 Use cases";                      C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 15; B = "2.4 Software Design / Components";   C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 16; B = "2.5 Interface Design";                C = 3;  D = 4;  E = 3;  F = 4;  G = 1 }
    @{ R = 17; B = "3.0 Implementation";                  C = 6;  D = 7;  E = $null; F = $null; G = 1 }
    @{ R = 18; B = "  3.1 Test Report";                    C = 6;  D = 7;  E = $null; F = $null; G = 1 }
    @{ R = 19; B = "3.2 Data processing";                  C = 5;  D = 3;  E = $null; F = $null; G = 1 }
    @{ R = 20; B = "3.3 User defined data recall";         C = 7;  D = 4;  E = $null; F = $null; G = 1 }
    @{ R = 21; B = "3.4 Data visualisation";               C = 7;  D = 4;  E = $null; F = $null; G = 1 }
    @{ R = 22; B = "3.5 Gui implementation";               C = 9;  D = 3;  E = $null; F = $null; G = 1 }
    @{ R = 23; B = "4.0 Controlling";                      C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 24; B = "4.1 Progress Meetings";                C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 25; B = "4.2 Status reports";                   C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 26; B = "4.3 Update plans";                     C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 27; B = "4.4 Update schedual";                  C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 28; B = "4.5 Version control";                  C = 1;  D = 12; E = $null; F = $null; G = $null }
    @{ R = 29; B = "5.0 Closing";                          C = 10; D = 3;  E = $null; F = $null; G = 1 }
    @{ R = 30; B = "5.1 Executive Summary";                C = 10; D = 3;  E = $null; F = $null; G = 1 }
    @{ R = 31; B = "5.2 User Manual";                      C = 10; D = 3;  E = $null; F = $null; G = 1 }
    @{ R = 32; B = "5.3 Self Assessment of Schedule";      C = 12; D = 1;  E = $null; F = $null; G = 1 }
    @{ R = 33; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 34; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 35; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 36; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 37; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 38; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
    @{ R = 39; B = $null; C = $null; D = $null; E = $null; F = $null; G = $null }
)

foreach ($row in $rows) {
    $r = $row.R

    if ($null -eq $row.B) {
        $ws.Range("B$r").ClearContents()
    } else {
        $ws.Range("B$r").Value = $row.B
    }

    foreach ($col in @("C", "D", "E", "F", "G")) {
        $val = $row[$col]
        if ($null -eq $val) {
            $ws.Range("$col$r").ClearContents()
        } else {
            $ws.Range("$col$r").Value = $val
        }
    }
}

# --- Title (set after the activity rows so the new shared string is
#     appended at the end of the table, matching the target layout) ----
$ws.Range("B1").Value = "Sydney Airbnb Data"

# --- Selected cell shown when the workbook is opened -------------------
[void]$ws.Range("D5").Select()

# --- Print scale 51% -> 50% --------------------------------------------
$ws.PageSetup.Zoom = 50
